$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-missing Start Time / End Time for row 51 (keep the
# same number formats as the cells above by copying first, then overwriting
# the value).
$ws.Range("B50").Copy($ws.Range("B51"))
$ws.Range("B51").Value = 0
$ws.Range("C50").Copy($ws.Range("C51"))
$ws.Range("C51").Value = 0.3354166666666667

# Add a new daily power record row (row 52): Date plus the three calculated
# columns (Start/End time are left blank, same as a fresh day).
$ws.Range("A51").Copy($ws.Range("A52"))
$ws.Range("A52").Value = 43376

$ws.Range("D51").Copy($ws.Range("D52"))
$ws.Range("D52").Formula = "=(C52-B52)* 1440"

$ws.Range("E51").Copy($ws.Range("E52"))
$ws.Range("E52").Formula = "=IF(C52>B52, (C52-B52)*1440, (B52-C52)*1440)"

$ws.Range("F51").Copy($ws.Range("F52"))
$ws.Range("F52").Formula = "=ABS((C52-B52)*1440)"

# Extend the table range so the new row is part of comforter_cda_table.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F52"))

# Keep the active selection / scrolled view in sync with the new last row.
$ws.Range("B52").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 41
